$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D column (prices, stored as text in the workbook) to stay
# text instead of being auto-converted to numbers when we assign values
# that look numeric.
$ws.Range("D2:D50").NumberFormat = "@"

# Simple price / value updates (D column)
$ws.Range("D2").Value  = "276.89"
$ws.Range("D4").Value  = "6.230"
$ws.Range("D5").Value  = "0.06199"
$ws.Range("D6").Value  = "3.574"
$ws.Range("D9").Value  = "0.8221"
$ws.Range("D10").Value = "0.1649"
$ws.Range("D11").Value = "0.08245"
$ws.Range("D12").Value = "0.03501"
$ws.Range("D14").Value = "0.09127"
$ws.Range("D15").Value = "3.771"
$ws.Range("D16").Value = "0.001623"
$ws.Range("D17").Value = "0.04703"
$ws.Range("D18").Value = "0.006468"
$ws.Range("D19").Value = "0.006148"
$ws.Range("D20").Value = "0.001067"
$ws.Range("D22").Value = "3.772"
$ws.Range("D24").Value = "0.01387"
$ws.Range("D28").Value = "0.0002736"
$ws.Range("D40").Value = "0.04676"

# Row 7 / 8 swap (KuCoinToken <-> FTXToken) plus price changes
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "1.548"
$ws.Range("E7").Value = "6FTXTokenFTT"

$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "6.570"
$ws.Range("E8").Value = "7KuCoinTokenKCS"

# Row 41 update
$ws.Range("D41").Value = "0.007021"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

# Row 42 / 43 swap (BKEXToken <-> CEJI) plus price changes
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.004700"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "0.1104"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# Remaining value updates
$ws.Range("D44").Value = "0.01084"
$ws.Range("D45").Value = "0.00006420"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D47").Value = "0.8451"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("D48").Value = "0.001391"
$ws.Range("D49").Value = "0.00001900"
$ws.Range("D50").Value = "0.01240"
